$d = $word.ActiveDocument

# 1. Replace the "Many AI techniques are involved in this project." paragraph with the new
#    combined decentralized/centralized summary paragraph.
$old1 = "Many AI techniques are involved in this project. "
$new1 = "For decentralized methods, each agent make their own decision and there are many levels of communication. Level-1: self-interested; level-2: sharing useful intel; level-3: sharing stats and maybe form tactics accordingly, with each level more information is shared. For centralized methods, there is a leader in the team who has all the information and give order to others. "
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# 2. Replace the "For decentralized methods, our thought is there are many levels..." paragraph
#    with "Our thought is using decision tree..." (content previously in the centralized paragraph).
$old2 = "For decentralized methods, our thought is there are many levels of communication. Level-1 self-interested; level-2 sharing useful intel; level-3 sharing stats and maybe form tactics accordingly and so on with each level more information is shared."
$new2 = "Our thought is using decision tree, machine learning techniques, state machines, game tree. And of course, there is some necessary AI needed like path-finding, decision of fight or flight and so on."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# 3. Delete the now-redundant "For centralized methods, our thought is using decision tree..."
#    paragraph entirely (its content was merged into the paragraph above). We remove the whole
#    paragraph, including its own paragraph mark, using a Range that spans from the start of
#    that paragraph to the start of the following (blank) paragraph, so the blank paragraph
#    that follows keeps its own (un-indented) formatting.
$target = "For centralized methods, our thought is using decision tree, machine learning techniques, state machines, game tree. And of course, there is some necessary AI needed like path-finding, decision of fight or flight and so on."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq $target) {
        $nextStart = $p.Next().Range.Start
        $r = $d.Range($p.Range.Start, $nextStart)
        $r.Delete()
        break
    }
}

# 4. Small wording fix in the Evaluation Method section: "to ask" -> "inveting".
$old4 = "Another way to evaluate would be to ask people to play against AI,"
$new4 = "Another way to evaluate would be inveting people to play against AI,"
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
